$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("J2").Value = 1.08
$ws.Range("L2").Value = 1.4
$ws.Range("O2").Value = 1.62

# Row 3 updates
$ws.Range("H3").Value = 2.75
$ws.Range("I3").Value = 2.7
$ws.Range("O3").Value = 1.6
$ws.Range("R3").Value = 1.83
$ws.Range("S3").Value = 1.83
